# Scheduled data refresh: update cached Universalis price snapshots and
# recomputed Leve-profit figures across the per-job "Behemoth_Profits" sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 477.18518
$ws.Range("I28").Value = 322.90475
$ws.Range("K28").Value = 322.90475
$ws.Range("M28").Value = 162.09525
$ws.Range("H116").Value = 5644.125
$ws.Range("I116").Value = 5108.154
$ws.Range("J116").Value = 7966.6665
$ws.Range("K116").Value = 5108.154
$ws.Range("L116").Value = 7966.6665
$ws.Range("M116").Value = -1666.154
$ws.Range("N116").Value = -14850.6665
$ws.Range("H132").Value = 1355.381
$ws.Range("I132").Value = 1198.15
$ws.Range("K132").Value = 3594.45
$ws.Range("M132").Value = -1064.45

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9095483
$ws.Range("I32").Value = 11365012
$ws.Range("K32").Value = 11365012
$ws.Range("M32").Value = -11364725
$ws.Range("H61").Value = 18300330
$ws.Range("I61").Value = 14712633
$ws.Range("J61").Value = 35726290
$ws.Range("K61").Value = 14712633
$ws.Range("L61").Value = 35726290
$ws.Range("M61").Value = -14712421
$ws.Range("N61").Value = -35726714
$ws.Range("H74").Value = 15300230
$ws.Range("I74").Value = 22730294
$ws.Range("J74").Value = 1678445.1
$ws.Range("K74").Value = 22730294
$ws.Range("L74").Value = 1678445.1
$ws.Range("M74").Value = -22729420
$ws.Range("N74").Value = -1680193.1
$ws.Range("H77").Value = 15300230
$ws.Range("I77").Value = 22730294
$ws.Range("J77").Value = 1678445.1
$ws.Range("K77").Value = 113651470
$ws.Range("L77").Value = 8392225.5
$ws.Range("M77").Value = -113647102
$ws.Range("N77").Value = -8400961.5
$ws.Range("H109").Value = 49461
$ws.Range("J109").Value = 49461
$ws.Range("L109").Value = 49461
$ws.Range("N109").Value = -52235
$ws.Range("H117").Value = 93250
$ws.Range("J117").Value = 93250
$ws.Range("L117").Value = 93250
$ws.Range("N117").Value = -102428
$ws.Range("H122").Value = 2937.7
$ws.Range("I122").Value = 2054.4285
$ws.Range("K122").Value = 6163.2855
$ws.Range("M122").Value = -3713.2855
$ws.Range("H136").Value = 18300330
$ws.Range("I136").Value = 14712633
$ws.Range("J136").Value = 35726290
$ws.Range("K136").Value = 44137899
$ws.Range("L136").Value = 107178870
$ws.Range("M136").Value = -44135349
$ws.Range("N136").Value = -107183970

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3269.8667
$ws.Range("I20").Value = 3580.6924
$ws.Range("J20").Value = 1249.5
$ws.Range("K20").Value = 3580.6924
$ws.Range("L20").Value = 1249.5
$ws.Range("M20").Value = -3333.6924
$ws.Range("N20").Value = -1743.5
$ws.Range("H107").Value = 2256.9
$ws.Range("I107").Value = 2333.625
$ws.Range("J107").Value = 1950
$ws.Range("K107").Value = 2333.625
$ws.Range("L107").Value = 1950
$ws.Range("M107").Value = -413.625
$ws.Range("N107").Value = -5790

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2073.75
$ws.Range("I58").Value = 1752.625
$ws.Range("J58").Value = 4000.5
$ws.Range("K58").Value = 1752.625
$ws.Range("L58").Value = 4000.5
$ws.Range("M58").Value = -1549.625
$ws.Range("N58").Value = -4406.5
$ws.Range("H86").Value = 5735.278
$ws.Range("I86").Value = 5340.636
$ws.Range("J86").Value = 6355.4287
$ws.Range("K86").Value = 5340.636
$ws.Range("L86").Value = 6355.4287
$ws.Range("M86").Value = -4217.636
$ws.Range("N86").Value = -8601.4287
$ws.Range("H89").Value = 5735.278
$ws.Range("I89").Value = 5340.636
$ws.Range("J89").Value = 6355.4287
$ws.Range("K89").Value = 26703.18
$ws.Range("L89").Value = 31777.1435
$ws.Range("M89").Value = -21087.18
$ws.Range("N89").Value = -43009.14350000001
$ws.Range("H136").Value = 2073.75
$ws.Range("I136").Value = 1752.625
$ws.Range("J136").Value = 4000.5
$ws.Range("K136").Value = 5257.875
$ws.Range("L136").Value = 12001.5
$ws.Range("M136").Value = -2707.875
$ws.Range("N136").Value = -17101.5
$ws.Range("H140").Value = 44099.5
$ws.Range("J140").Value = 44099
$ws.Range("L140").Value = 44099
$ws.Range("N140").Value = -54459

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 45000
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").Value = ""
$ws.Range("H68").Value = 1749.3334
$ws.Range("J68").Value = 2900
$ws.Range("L68").Value = 8700
$ws.Range("N68").Value = -10322
$ws.Range("H71").Value = 1749.3334
$ws.Range("J71").Value = 2900
$ws.Range("L71").Value = 26100
$ws.Range("N71").Value = -34212
$ws.Range("H81").Value = 500
$ws.Range("J81").Value = 500
$ws.Range("L81").Value = 1500
$ws.Range("N81").Value = -3746
$ws.Range("H84").Value = 500
$ws.Range("J84").Value = 500
$ws.Range("L84").Value = 4500
$ws.Range("N84").Value = -15732
$ws.Range("H107").Value = 593.9231
$ws.Range("J107").Value = 792.1111
$ws.Range("L107").Value = 2376.3333
$ws.Range("N107").Value = -6216.3333
$ws.Range("H131").Value = 6843.1294
$ws.Range("J131").Value = 6966.1133
$ws.Range("L131").Value = 20898.3399
$ws.Range("N131").Value = -30978.3399
$ws.Range("H139").Value = 3656.963
$ws.Range("I139").Value = 1978.8334
$ws.Range("K139").Value = 5936.5002
$ws.Range("M139").Value = -796.5002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 25000
$ws.Range("J48").Value = 25000
$ws.Range("L48").Value = 25000
$ws.Range("N48").Value = -25970
$ws.Range("H70").Value = 11180.625
$ws.Range("I70").Value = 15861.75
$ws.Range("K70").Value = 15861.75
$ws.Range("M70").Value = -15591.75
$ws.Range("H73").Value = 11180.625
$ws.Range("I73").Value = 15861.75
$ws.Range("K73").Value = 15861.75
$ws.Range("M73").Value = -14925.75
$ws.Range("H122").Value = 1592.5
$ws.Range("I122").Value = 1688.8889
$ws.Range("J122").Value = 1303.3334
$ws.Range("K122").Value = 5066.6667
$ws.Range("L122").Value = 3910.0002
$ws.Range("M122").Value = -2616.6667
$ws.Range("N122").Value = -8810.0002
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = ""
$ws.Range("N123").Value = ""
$ws.Range("H132").Value = 32265852
$ws.Range("I132").Value = 58827468
$ws.Range("J132").Value = 12460.643
$ws.Range("K132").Value = 176482404
$ws.Range("L132").Value = 37381.929
$ws.Range("M132").Value = -176479874
$ws.Range("N132").Value = -42441.929
$ws.Range("H136").Value = 10113
$ws.Range("J136").Value = 10113
$ws.Range("L136").Value = 30339
$ws.Range("N136").Value = -35439

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1459.5333
$ws.Range("I22").Value = 1542.7778
$ws.Range("J22").Value = 1334.6666
$ws.Range("K22").Value = 1542.7778
$ws.Range("L22").Value = 1334.6666
$ws.Range("M22").Value = -1247.7778
$ws.Range("N22").Value = -1924.6666
$ws.Range("H27").Value = 1459.5333
$ws.Range("I27").Value = 1542.7778
$ws.Range("J27").Value = 1334.6666
$ws.Range("K27").Value = 1542.7778
$ws.Range("L27").Value = 1334.6666
$ws.Range("M27").Value = -1435.7778
$ws.Range("N27").Value = -1548.6666
$ws.Range("H55").Value = 52632330
$ws.Range("I55").Value = 66667376
$ws.Range("K55").Value = 66667376
$ws.Range("M55").Value = -66667203
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = ""
$ws.Range("N80").Value = ""
$ws.Range("H81").Value = 69909
$ws.Range("J81").Value = 69909
$ws.Range("L81").Value = 69909
$ws.Range("N81").Value = -71905
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = ""
$ws.Range("N83").Value = ""
$ws.Range("H84").Value = 69909
$ws.Range("J84").Value = 69909
$ws.Range("L84").Value = 209727
$ws.Range("N84").Value = -219711
$ws.Range("H100").Value = 4490.6665
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 4490.6665
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = ""
$ws.Range("M100").Value = ""
$ws.Range("N100").Value = -5572.6665
$ws.Range("H110").Value = 48199.8
$ws.Range("J110").Value = 48199.8
$ws.Range("L110").Value = 48199.8
$ws.Range("N110").Value = -56379.8
$ws.Range("H132").Value = 1133645.5
$ws.Range("I132").Value = 59434.668
$ws.Range("J132").Value = 1670750.9
$ws.Range("K132").Value = 178304.004
$ws.Range("L132").Value = 5012252.699999999
$ws.Range("M132").Value = -175774.004
$ws.Range("N132").Value = -5017312.699999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 23831.666
$ws.Range("I49").Value = 19000
$ws.Range("J49").Value = 33495
$ws.Range("K49").Value = 19000
$ws.Range("L49").Value = 33495
$ws.Range("M49").Value = -18770
$ws.Range("N49").Value = -33955
$ws.Range("H107").Value = 849.19354
$ws.Range("I107").Value = 905.5217
$ws.Range("J107").Value = 687.25
$ws.Range("K107").Value = 2716.5651
$ws.Range("L107").Value = 2061.75
$ws.Range("M107").Value = -796.5650999999998
$ws.Range("N107").Value = -5901.75
$ws.Range("H122").Value = 2346.739
$ws.Range("I122").Value = 2289.2856
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 6867.8568
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -4417.8568
$ws.Range("N122").Value = -13750
